$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17/18 swap: WrappedBTC <-> Chainlink ---
$ws.Range("B17").Value2 = "Chainlink"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "18.10"
$ws.Range("E17").Value2 = "  -2.27%  "
$ws.Range("B18").Value2 = "WrappedBTC"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "66.198.53"
$ws.Range("E18").Value2 = "  -4.84%  "

# --- Row 24/25 swap: Toncoin <-> Litecoin ---
$ws.Range("B24").Value2 = "Litecoin"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "101.04"
$ws.Range("E24").Value2 = "  -4.32%  "
$ws.Range("B25").Value2 = "Toncoin"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "5.03"
$ws.Range("E25").Value2 = "  -4.59%  "

# --- Remaining price/volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "66.226.10"
$ws.Range("E2").Value2 = "  -4.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.348.61"
$ws.Range("E3").Value2 = "  -5.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("E4").Value2 = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "561.72"
$ws.Range("E5").Value2 = "  -3.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "182.50"
$ws.Range("E6").Value2 = "  -7.74%  "
$ws.Range("E7").Value2 = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.592"
$ws.Range("E8").Value2 = "  -3.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "3.343.42"
$ws.Range("E9").Value2 = "  -5.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.185"
$ws.Range("E10").Value2 = "  -9.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.589"
$ws.Range("E11").Value2 = "  -6.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "47.51"
$ws.Range("E12").Value2 = "  -8.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000265"
$ws.Range("E13").Value2 = "  -7.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "8.65"
$ws.Range("E14").Value2 = "  -6.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "3.872.86"
$ws.Range("E15").Value2 = "  -5.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "605.94"
$ws.Range("E16").Value2 = "  -8.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "3.339.76"
$ws.Range("E19").Value2 = "  -5.64%  "
$ws.Range("E20").Value2 = "  -3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "11.45"
$ws.Range("E21").Value2 = "  -8.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.908"
$ws.Range("E22").Value2 = "  -6.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "16.90"
$ws.Range("E23").Value2 = "  -7.94%  "
$ws.Range("E26").Value2 = "  -7.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "6.01"
$ws.Range("E27").Value2 = "  -0.45%  "
$ws.Range("E28").Value2 = "  -8.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "9.32"
$ws.Range("E29").Value2 = "  -8.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "8.73"
$ws.Range("E30").Value2 = "  -9.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "30.44"
$ws.Range("E31").Value2 = "  -9.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "6.27"
$ws.Range("E32").Value2 = "  -7.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.80"
$ws.Range("E33").Value2 = "  -14.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "11.06"
$ws.Range("E34").Value2 = "  -6.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.835.57"
$ws.Range("E35").Value2 = "  +1.26%  "
$ws.Range("E36").Value2 = "  -5.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "532.64"
$ws.Range("E37").Value2 = "  +5.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "57.57"
$ws.Range("E39").Value2 = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "3.41"
$ws.Range("E40").Value2 = "  -9.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0₃0715"
$ws.Range("E41").Value2 = "  -12.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "2.66"
$ws.Range("E42").Value2 = "  -9.17%  "
$ws.Range("E43").Value2 = "  -7.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.343"
$ws.Range("E44").Value2 = "  -8.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "32.07"
$ws.Range("E45").Value2 = "  -7.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.19"
$ws.Range("E46").Value2 = "  +18.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.0413"
$ws.Range("E47").Value2 = "  -8.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "3.13"
$ws.Range("E48").Value2 = "  -7.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.130"
$ws.Range("E49").Value2 = "  -4.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.62"
$ws.Range("E50").Value2 = "  -9.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.999"
$ws.Range("E51").Value2 = "  -0.21%  "
